$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'66.524.81"
$ws.Range("E2").Value = "  -0.86%  "

# Row 3
$ws.Range("D3").Value = "'3.801.52"
$ws.Range("E3").Value = "  +0.43%  "

# Row 4
$ws.Range("D4").Value = "'0.997"
$ws.Range("E4").Value = "  -0.35%  "

# Row 5
$ws.Range("D5").Value = "'419.82"
$ws.Range("E5").Value = "  +0.96%  "

# Row 6
$ws.Range("D6").Value = "'128.58"
$ws.Range("E6").Value = "  -7.18%  "

# Row 7
$ws.Range("D7").Value = "'3.800.39"
$ws.Range("E7").Value = "  +0.55%  "

# Row 8
$ws.Range("D8").Value = "'0.599"
$ws.Range("E8").Value = "  -6.18%  "

# Row 9
$ws.Range("E9").Value = "  +0.10%  "

# Row 10
$ws.Range("D10").Value = "'0.718"
$ws.Range("E10").Value = "  -5.21%  "

# Row 11
$ws.Range("D11").Value = "'0.161"
$ws.Range("E11").Value = "  -9.28%  "

# Row 12
$ws.Range("D12").Value = "'0.0000345"
$ws.Range("E12").Value = "  -7.04%  "

# Row 13
$ws.Range("D13").Value = "'39.75"
$ws.Range("E13").Value = "  -6.41%  "

# Row 14
$ws.Range("D14").Value = "'4.377.34"
$ws.Range("E14").Value = "  -0.01%  "

# Row 15
$ws.Range("E15").Value = "  -2.45%  "

# Row 16
$ws.Range("D16").Value = "'15.90"
$ws.Range("E16").Value = "  +19.33%  "

# Row 17
$ws.Range("E17").Value = "  -1.65%  "

# Row 18
$ws.Range("D18").Value = "'3.792.04"
$ws.Range("E18").Value = "  -0.40%  "

# Row 19
$ws.Range("D19").Value = "'19.37"
$ws.Range("E19").Value = "  -5.47%  "

# Row 20
$ws.Range("D20").Value = "'66.546.21"
$ws.Range("E20").Value = "  -1.06%  "

# Row 21
$ws.Range("E21").Value = "  -3.89%  "

# Row 22
$ws.Range("D22").Value = "'404.00"
$ws.Range("E22").Value = "  -7.41%  "

# Row 23
$ws.Range("D23").Value = "'14.16"
$ws.Range("E23").Value = "  -6.59%  "

# Row 24
$ws.Range("D24").Value = "'83.30"
$ws.Range("E24").Value = "  -6.36%  "

# Row 25
$ws.Range("D25").Value = "'2.97"
$ws.Range("E25").Value = "  -2.85%  "

# Row 26
$ws.Range("D26").Value = "'36.79"
$ws.Range("E26").Value = "  -1.51%  "

# Row 27
$ws.Range("E27").Value = "  +10.40%  "

# Row 28
$ws.Range("D28").Value = "'3.18"
$ws.Range("E28").Value = "  -2.62%  "

# Row 29
$ws.Range("D29").Value = "'9.28"
$ws.Range("E29").Value = "  -4.74%  "

# Row 30
$ws.Range("D30").Value = "'700.01"
$ws.Range("E30").Value = "  -1.12%  "

# Row 31
$ws.Range("D31").Value = "'8.14"
$ws.Range("E31").Value = "  +13.74%  "

# Row 32
$ws.Range("D32").Value = "'2.77"
$ws.Range("E32").Value = "  +1.37%  "

# Row 33
$ws.Range("D33").Value = "'12.21"
$ws.Range("E33").Value = "  -2.54%  "

# Row 34
$ws.Range("E34").Value = "  -1.95%  "

# Row 35
$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.11%  "

# Row 36
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.150"
$ws.Range("E36").Value = "  -6.81%  "

# Row 37
$ws.Range("D37").Value = "'37.80"
$ws.Range("E37").Value = "  -8.51%  "

# Row 38
$ws.Range("D38").Value = "'54.48"
$ws.Range("E38").Value = "  -5.58%  "

# Row 39
$ws.Range("D39").Value = "'0.0₃0768"
$ws.Range("E39").Value = "  +13.31%  "

# Row 40
$ws.Range("E40").Value = "  -6.74%  "

# Row 41
$ws.Range("D41").Value = "'2.96"
$ws.Range("E41").Value = "  -0.90%  "

# Row 42
$ws.Range("D42").Value = "'4.67"
$ws.Range("E42").Value = "  +7.57%  "

# Row 43
$ws.Range("D43").Value = "'0.999"
$ws.Range("E43").Value = "  +0.06%  "

# Row 44
$ws.Range("E44").Value = "  -7.13%  "

# Row 45
$ws.Range("D45").Value = "'3.31"
$ws.Range("E45").Value = "  -1.79%  "

# Row 46
$ws.Range("D46").Value = "'144.49"
$ws.Range("E46").Value = "  -1.92%  "

# Row 47
$ws.Range("D47").Value = "'3.06"
$ws.Range("E47").Value = "  -1.18%  "

# Row 48
$ws.Range("E48").Value = "  -2.70%  "

# Row 49
$ws.Range("D49").Value = "'25.53"
$ws.Range("E49").Value = "  -3.63%  "

# Row 50
$ws.Range("B50").Value = "WEMIXToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").Value = "'2.51"
$ws.Range("E50").Value = "  -2.31%  "

# Row 51
$ws.Range("B51").Value = "Stacks"
$ws.Range("C51").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D51").Value = "'2.73"
$ws.Range("E51").Value = "  -3.91%  "
